$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values remain stored as text, matching the
# original inlineStr cell type, instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.438.17'
$ws.Range("E2").Value = '  -3.85%  '
$ws.Range("D3").Value = '1.772.58'
$ws.Range("E3").Value = '  -2.90%  '
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '306.58'
$ws.Range("E6").Value = '  -1.85%  '
$ws.Range("D7").Value = '0.4298'
$ws.Range("E7").Value = '  +1.07%  '
$ws.Range("D8").Value = '0.3673'
$ws.Range("E8").Value = '  +1.69%  '
$ws.Range("D9").Value = '0.07231'
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '0.8484'
$ws.Range("E10").Value = '  -1.57%  '
$ws.Range("D11").Value = '20.32'
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").Value = '1.775.70'
$ws.Range("E12").Value = '  -3.63%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '6.440'
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '5.248'
$ws.Range("E14").Value = '  -2.69%  '
$ws.Range("D15").Value = '0.06854'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").Value = '79.61'
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").Value = '0.000008689'
$ws.Range("E18").Value = '  -2.61%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '15.06'
$ws.Range("E20").Value = '  -1.97%  '
$ws.Range("D21").Value = '26.440.14'
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("D22").Value = '5.108'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '11.29'
$ws.Range("E23").Value = '  +3.65%  '
$ws.Range("D24").Value = '1.997.10'
$ws.Range("E24").Value = '  -2.99%  '
$ws.Range("D25").Value = '152.39'
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("E26").Value = '  -6.97%  '
$ws.Range("D27").Value = '18.16'
$ws.Range("E27").Value = '  -2.93%  '
$ws.Range("D28").Value = '5.098'
$ws.Range("D29").Value = '114.89'
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("D30").Value = '1.718'
$ws.Range("E30").Value = '  -4.25%  '
$ws.Range("D32").Value = '0.7242'
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("D33").Value = '1.118'
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").Value = '4.335'
$ws.Range("E34").Value = '  -4.56%  '
$ws.Range("D35").Value = '2.755'
$ws.Range("E35").Value = '  -7.29%  '
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("D37").Value = '1.078'
$ws.Range("E37").Value = '  -0.92%  '
$ws.Range("D38").Value = '0.05161'
$ws.Range("E39").Value = '  -1.29%  '
$ws.Range("D40").Value = '0.4928'
$ws.Range("E40").Value = '  -2.90%  '
$ws.Range("E41").Value = '  -2.96%  '
$ws.Range("D42").Value = '2.533'
$ws.Range("E42").Value = '  -9.32%  '
$ws.Range("D43").Value = '6.215'
$ws.Range("E43").Value = '  -3.76%  '
$ws.Range("D44").Value = '8.042'
$ws.Range("E44").Value = '  -3.78%  '
$ws.Range("D45").Value = '104.90'
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").Value = '1.002'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").Value = '10.17'
$ws.Range("E47").Value = '  -3.07%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.06196'
$ws.Range("E48").Value = '  -4.20%  '
$ws.Range("B49").Value = 'Decentraland'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D49").Value = '0.4488'
$ws.Range("E49").Value = '  -4.04%  '
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").Value = '1.743'
$ws.Range("E51").Value = '  +1.56%  '
